$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")
$summary = $wb.Worksheets.Item("Summary")

function Set-TextValue {
    param($Cell, $Text)
    # This workbook stores every value (including numeric-looking ones like
    # "3" or "15") as text - the sheet even carries a
    # numberStoredAsText ignoredError for the whole range. Forcing the
    # NumberFormat to "@" before assigning the value keeps new cells text
    # too, and resetting the style back to "Normal" afterwards avoids
    # leaving a stray number-format override on the cell.
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

# Order 3 and order 4 line items appended below the existing rows 1-11
# (new used range becomes A1:L21).
$newRows = @(
    @{ Row = 12; A = "3"; C = "192_粉荔枝_Pink Ohara_Rosa rugosa Thunb._20stems"; F = "3" },
    @{ Row = 13;          C = "197_粉红雪山_Sweet Avalanche_Rosa rugosa Thunb._20stems"; F = "15" },
    @{ Row = 14;          C = "192_粉荔枝_Pink Ohara_Rosa rugosa Thunb._20stems"; F = "10" },
    @{ Row = 15;          C = "148_坦尼克_Tineke_Rosa rugosa Thunb._20stems"; F = "19" },
    @{ Row = 16;          C = "479_绿灵草_lepidium_undefined_1bunch"; F = "25" },
    @{ Row = 17; A = "4"; C = "138_卡罗拉_Carola_Rosa rugosa Thunb._20stems"; F = "18" },
    @{ Row = 18;          C = "148_坦尼克_Tineke_Rosa rugosa Thunb._20stems"; F = "4" },
    @{ Row = 19;          C = "173_朱丽叶_Juliet_Rosa rugosa Thunb._20stems"; F = "11" },
    @{ Row = 20;          C = "209_海洋之歌_Ocean Song_Rosa rugosa Thunb._20stems"; F = "8" },
    @{ Row = 21;          C = "184_微光_shimmer_Rosa rugosa Thunb._20stems" }
)

foreach ($item in $newRows) {
    $r = $item.Row
    if ($item.ContainsKey("A")) {
        Set-TextValue $ws.Cells.Item($r, 1) $item.A
    }
    $ws.Cells.Item($r, 3).Value = $item.C
    if ($item.ContainsKey("F")) {
        Set-TextValue $ws.Cells.Item($r, 6) $item.F
    }
}

# Summary sheet: the tracking/awb number in G2 grew additional digits.
Set-TextValue $summary.Range("G2") "01013673102815383151019251841180"
